# Generate Report for Handoff
# Moves the localization status from "In Translation" to "Ready for handoff"
# and refreshes the related timestamps on all three sheets.

$wb = $excel.ActiveWorkbook

# --- Overview sheet ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsOverview.Range("G2").Value = "2016-09-01 09:15:43"
$wsOverview.Columns.Item(5).ColumnWidth = 16.3
$wsOverview.Columns.Item(6).ColumnWidth = 16.3

# --- zh-cn sheet ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsZhCn.Range("H2").Value = "2016-09-01 09:15:39"
$wsZhCn.Columns.Item(3).ColumnWidth = 16.3

# --- de-de sheet ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("H2").Value = "2016-09-01 09:15:43"
$wsDeDe.Columns.Item(3).ColumnWidth = 16.3
